# "resolved security group bug" -- every shape in the architecture
# diagram on slide 1 shifts by the same offset (dx=-150019 EMU,
# dy=-196931 EMU), so the group lines back up correctly.
#
# Shape.Left/Top on this host are 32-bit (Single) point values that get
# multiplied by 12700 and truncated back to EMU on save, so each target
# below is the nearest point value whose Single representation reproduces
# the exact target EMU from the diff (verified against run_com output).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# "Rechteck 21" (id 22): off (2886068,5689814) -> (2736049,5492883) EMU
$shp = $s.Shapes.Item("Rechteck 21")
$shp.Left = 215.4369354248047
$shp.Top = 432.510498046875

# "Rechteck 22" (id 23): off (2897976,4612366) -> (2747957,4415435) EMU
$shp = $s.Shapes.Item("Rechteck 22")
$shp.Left = 216.37457275390625
$shp.Top = 347.67205810546875

# "Rechteck 3" (id 4): off (3121814,5819775) -> (2971795,5622844) EMU
$shp = $s.Shapes.Item("Rechteck 3")
$shp.Left = 233.99961853027344
$shp.Top = 442.74365234375

# "Rechteck 7" (id 8): off (4431504,4757737) -> (4281485,4560806) EMU
$shp = $s.Shapes.Item("Rechteck 7")
$shp.Left = 337.12481689453125
$shp.Top = 359.11859130859375

# "Rechteck 20" (id 21): off (2897976,3644502) -> (2747957,3447571) EMU
$shp = $s.Shapes.Item("Rechteck 20")
$shp.Left = 216.37457275390625
$shp.Top = 271.4623107910156

# "Rechteck 8" (id 9): off (3121813,3695701) -> (2971794,3498770) EMU
$shp = $s.Shapes.Item("Rechteck 8")
$shp.Left = 233.99954223632812
$shp.Top = 275.49371337890625

# "Rechteck 9" (id 10): off (3121815,4757738) -> (2971796,4560807) EMU
$shp = $s.Shapes.Item("Rechteck 9")
$shp.Left = 233.99969482421875
$shp.Top = 359.1186828613281

# "Rechteck 10" (id 11): off (4431504,3695701) -> (4281485,3498770) EMU
$shp = $s.Shapes.Item("Rechteck 10")
$shp.Left = 337.12481689453125
$shp.Top = 275.49371337890625

# "Rechteck 11" (id 12): off (6927060,3695700) -> (6777041,3498769) EMU
$shp = $s.Shapes.Item("Rechteck 11")
$shp.Left = 533.6253051757812
$shp.Top = 275.49365234375

# "Rechteck 12" (id 13): off (8222461,3695700) -> (8072442,3498769) EMU
$shp = $s.Shapes.Item("Rechteck 12")
$shp.Left = 635.6253662109375
$shp.Top = 275.49365234375

# "Rechteck 13" (id 14): off (6927060,4757737) -> (6777041,4560806) EMU
$shp = $s.Shapes.Item("Rechteck 13")
$shp.Left = 533.6253051757812
$shp.Top = 359.11859130859375

# "Rechteck 14" (id 15): off (8222462,4757736) -> (8072443,4560805) EMU
$shp = $s.Shapes.Item("Rechteck 14")
$shp.Left = 635.62548828125
$shp.Top = 359.1185302734375

# "Rechteck 15" (id 16): off (2774156,3557588) -> (2624137,3360657) EMU
$shp = $s.Shapes.Item("Rechteck 15")
$shp.Left = 206.62496948242188
$shp.Top = 264.6186828613281

# "Rechteck 16" (id 17): off (6573447,3557588) -> (6423428,3360657) EMU
$shp = $s.Shapes.Item("Rechteck 16")
$shp.Left = 505.78173828125
$shp.Top = 264.6186828613281

# "Textfeld 17" (id 18): off (4276722,6242024) -> (4126703,6045093) EMU
$shp = $s.Shapes.Item("Textfeld 17")
$shp.Left = 324.937255859375
$shp.Top = 475.9915771484375

# "Textfeld 19" (id 20): off (7998623,6228874) -> (7848604,6031943) EMU
$shp = $s.Shapes.Item("Textfeld 19")
$shp.Left = 618.0003662109375
$shp.Top = 474.9561462402344

# "Textfeld 23" (id 24): off (10042924,3871788) -> (9892905,3674857) EMU
$shp = $s.Shapes.Item("Textfeld 23")
$shp.Left = 778.9689331054688
$shp.Top = 289.35882568359375

# "Textfeld 25" (id 26): off (10042923,4915972) -> (9892904,4719041) EMU
$shp = $s.Shapes.Item("Textfeld 25")
$shp.Left = 778.9688720703125
$shp.Top = 371.5780334472656

# "Textfeld 26" (id 27): off (10036970,5859542) -> (9886951,5662611) EMU
$shp = $s.Shapes.Item("Textfeld 26")
$shp.Left = 778.5001220703125
$shp.Top = 445.8749084472656

# "Rechteck 27" (id 28): off (2150269,3429000) -> (2000250,3232069) EMU
$shp = $s.Shapes.Item("Rechteck 27")
$shp.Left = 157.5
$shp.Top = 254.49363708496094

# "Textfeld 28" (id 29): off (2228400,3943350) -> (2078381,3746419) EMU
$shp = $s.Shapes.Item("Textfeld 28")
$shp.Left = 163.6520538330078
$shp.Top = 294.99365234375

# "Rechteck 30" (id 31): off (4772845,1703294) -> (4622826,1506363) EMU
$shp = $s.Shapes.Item("Rechteck 30")
$shp.Left = 364.0020751953125
$shp.Top = 118.61126708984375

# "Gerade Verbindung mit Pfeil 32" (id 33): off (4358878,2348179) -> (4208859,2151248) EMU
$shp = $s.Shapes.Item("Gerade Verbindung mit Pfeil 32")
$shp.Left = 331.40625
$shp.Top = 169.38961791992188

# "Gerade Verbindung mit Pfeil 33" (id 34): off (6179346,2348179) -> (6029327,2151248) EMU
$shp = $s.Shapes.Item("Gerade Verbindung mit Pfeil 33")
$shp.Left = 474.75018310546875
$shp.Top = 169.38961791992188

# "Grafik 37" (id 38): off (5722146,96918) -> (5572127,-100013) EMU
$shp = $s.Shapes.Item("Grafik 37")
$shp.Left = 438.75018310546875
$shp.Top = -7.875039577484131

# "Gerade Verbindung mit Pfeil 39" (id 40): off (6179346,1011318) -> (6029327,814387) EMU
$shp = $s.Shapes.Item("Gerade Verbindung mit Pfeil 39")
$shp.Left = 474.75018310546875
$shp.Top = 64.12496185302734

# "Rechteck 40" (id 41): off (1343023,1516556) -> (1193004,1319625) EMU
$shp = $s.Shapes.Item("Rechteck 40")
$shp.Left = 93.93732452392578
$shp.Top = 103.90748596191406

# "Textfeld 42" (id 43): off (1407288,1590957) -> (1257269,1394026) EMU
$shp = $s.Shapes.Item("Textfeld 42")
$shp.Left = 98.99756622314453
$shp.Top = 109.76583099365234
